# Correção das notas do fórum para matc65 em 2021.2
# For every student row where the "nota_view" column (J) equals 4,
# reset the daily-view flags and totals (columns B..J) to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($row = 2; $row -le $lastRow; $row++) {
    $nota = $ws.Cells.Item($row, 10).Value2  # Column J = nota_view
    if ($nota -eq 4) {
        for ($col = 2; $col -le 10; $col++) {  # Columns B..J
            $ws.Cells.Item($row, $col).Value = 0
        }
    }
}
